$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived values for rows 2-5, columns E-H and M-T

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.370913
$ws.Range("H2").Value = 1.112739
$ws.Range("M2").Value = 5.740110333333334
$ws.Range("N2").Value = 17.220331
$ws.Range("O2").Value = 0.2861925343043439
$ws.Range("P2").Value = 0.2861925343043439
$ws.Range("Q2").Value = 2.129081544067667
$ws.Range("R2").Value = 19.161733896609
$ws.Range("S2").Value = 0.2861925343043439
$ws.Range("T2").Value = 0.2861925343043439

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.370913
$ws.Range("H3").Value = 1.112739
$ws.Range("O3").Value = 0.2917347240316885
$ws.Range("P3").Value = 0.2917347240316885
$ws.Range("Q3").Value = 2.170311738597
$ws.Range("R3").Value = 19.532805647373
$ws.Range("S3").Value = 0.2917347240316885
$ws.Range("T3").Value = 0.2917347240316885

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.370913
$ws.Range("H4").Value = 1.112739
$ws.Range("M4").Value = 6.759986
$ws.Range("N4").Value = 20.279958
$ws.Range("O4").Value = 0.3370418707750538
$ws.Range("P4").Value = 0.3370418707750538
$ws.Range("Q4").Value = 2.507366687218
$ws.Range("R4").Value = 22.566300184962
$ws.Range("S4").Value = 0.3370418707750538
$ws.Range("T4").Value = 0.3370418707750538

$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.370913
$ws.Range("H5").Value = 1.112739
$ws.Range("M5").Value = 1.705448333333333
$ws.Range("N5").Value = 5.116345
$ws.Range("O5").Value = 0.0850308708889137
$ws.Range("P5").Value = 0.0850308708889137
$ws.Range("Q5").Value = 0.6325729576616667
$ws.Range("R5").Value = 5.693156618954999
$ws.Range("S5").Value = 0.0850308708889137
$ws.Range("T5").Value = 0.0850308708889137
